# Updated symbol list on Tue Dec 20 18:23:27 UTC 2022 with GitHub Actions
#
# Refresh the "Price" (D) column with newly-scraped quotes, and fix up a
# couple of "Worstin24h"/"Bestin24h" suffix glitches in the Volume(1h)
# (E) column that the scraper mis-placed.
#
# NOTE: the Price column is stored as text (not numbers) in this sheet, so
# plain numeric literals that look like "249.07" must be written back as
# text -- otherwise Excel auto-converts them to floating point numbers and
# we lose exact formatting (trailing zeros, etc.) plus the cell's stored
# type changes. Forcing NumberFormat="@" (Text) before the write keeps the
# value a string; ClearFormats() afterwards drops the now-unneeded explicit
# number-format style so the cell's style index is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = [ordered]@{
    "D2"  = "249.07"
    "D3"  = "22.90"
    "D4"  = "5.425"
    "D5"  = "0.05621"
    "D6"  = "3.425"
    "D7"  = "6.358"
    "D8"  = "0.8144"
    "D9"  = "0.8973"
    "D10" = "0.1436"
    "D11" = "0.07494"
    "D12" = "0.03104"
    "D13" = "0.03101"
    "D14" = "0.09326"
    "D15" = "3.562"
    "D16" = "0.001580"
    "D17" = "0.04757"
    "D18" = "0.0005794"
    "D19" = "0.006380"
    "D20" = "0.004989"
    "D21" = "0.001033"
    "D22" = "0.0001502"
    "D23" = "3.698"
    "D24" = "2.191"
    "D25" = "0.3299"
    "D26" = "0.1296"
    "D28" = "0.0003034"
    "D40" = "0.04027"
    "D41" = "0.006848"
    "D42" = "0.1065"
    "D43" = "0.002719"
    "D44" = "0.007768"
    "D45" = "0.00005513"
    "D48" = "0.2405"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# Volume(1h) label glitches -- plain text, no numeric coercion risk.
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E27").Value = "26AAXTokenAAB"
